# position sync create a building by player's position
#
# Skill balance pass: tighten skill cooldown (column P, "CoolDownTime") and
# damage distance (column R, "DamageDistance") on the "Property1" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# --- Column R (DamageDistance): every data row from 11 to 61 goes from 2.5 to 2
for ($row = 11; $row -le 61; $row++) {
    $ws.Cells.Item($row, 18).Value = 2
}

# --- Column P (CoolDownTime): only the rows that currently hold the "short"
# cooldown value (2) move down to 1. These are rows 11-19 and the "Normal"
# skill row of every 3-row monster block from row 38 onward (38, 41, 47, 50,
# 53, 56, 59).
$pRows = @(11,12,13,14,15,16,17,18,19,38,41,47,50,53,56,59)
foreach ($row in $pRows) {
    $ws.Cells.Item($row, 16).Value = 1
}
